$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.153521299362183
$ws.Range("B1").Value = 5.941781997680664
$ws.Range("C1").Value = 1.237406492233276
$ws.Range("D1").Value = 0.246417224407196
$ws.Range("E1").Value = 0.3339578211307526
